# Append a new ledger entry (row 44) to the ledger sheet, mirroring the
# previous row (43), and update the summary formulas that referenced the old
# last row (43) so they now reference the new last row (44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 43 (values + formatting) down into row 44 so the new row keeps the
# exact same number formats / styles as the row above it.
$ws.Range("B43:M43").Copy($ws.Range("B44:M44"))

# New entry values for row 44 (mirrors the "initiation" status of row 43).
$ws.Range("B44").Value() = "2018.09.22 02:38:53"
$ws.Range("C44").Value() = "initiation"
$ws.Range("D44").Value() = 344379
$ws.Range("E44").Value() = 152.999236
$ws.Range("F44").Value() = 251829.180347625
$ws.Range("G44").Value() = 571.991631325
$ws.Range("H44").Value() = 596208.1803476249
$ws.Range("I44").Value() = 724.990867325

# Recreate the per-row formulas for row 44, referencing the new previous row.
$ws.Range("J44").Formula() = '=IF(C44="settlement", H44-H43, "")'
$ws.Range("K44").Formula() = '=IF(C44="settlement", I44-I43, "")'
$ws.Range("L44").Formula() = '=IF(C44="settlement", J44/H43, "")'
$ws.Range("M44").Formula() = '=IF(C44="settlement", SUM($J$11:J44)/$H$11, "")'

# Update the summary cells that used to point at the last row (43) so they
# now point at the new last row (44).
$ws.Range("C5").Formula() = '=SUM(J11:J44)'
$ws.Range("C6").Formula() = '=SUM(K11:K44)'
$ws.Range("C7").Formula() = '=M44'
